{"js": "// Fill in the next empty row of the team roster table with the new\n// member's Roll#, Name and Position title (team4 doc: \"Ryo Nishihira\").\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"values\");\nawait context.sync();\n\n// Locate the first fully-empty data row (skip the header row).\nconst rowValues = table.values;\nlet targetRow = -1;\nfor (let r = 1; r < rowValues.length; r++) {\n  if (rowValues[r].every((cell) => !cell || cell.trim() === \"\")) {\n    targetRow = r;\n    break;\n  }\n}\nif (targetRow === -1) {\n  targetRow = rowValues.length - 1;\n}\n\ntable.getCell(targetRow, 0).value = \"25\";\ntable.getCell(targetRow, 1).value = \"Ryo Nishihira\";\ntable.getCell(targetRow, 2).value = \"Developer\";\n\nawait context.sync();\n", "ps1": "# Fill in the next empty row of the team roster table with the new\n# member's Roll#, Name and Position title (team4 doc: \"Ryo Nishihira\").\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n# Find the first fully empty data row (skip the bold header row).\n$targetRow = -1\nfor ($r = 2; $r -le $t.Rows.Count; $r++) {\n    $allEmpty = $true\n    for ($c = 1; $c -le $t.Columns.Count; $c++) {\n        $cellText = $t.Cell($r, $c).Range.Text\n        $cellText = $cellText -replace \"[\\x07\\x0d]\", \"\"\n        if ($cellText.Trim() -ne \"\") {\n            $allEmpty = $false\n            break\n        }\n    }\n    if ($allEmpty) {\n        $targetRow = $r\n        break\n    }\n}\nif ($targetRow -eq -1) {\n    $targetRow = $t.Rows.Count\n}\n\n$values = @(\"25\", \"Ryo Nishihira\", \"Developer\")\nfor ($c = 1; $c -le 3; $c++) {\n    $cellRange = $t.Cell($targetRow, $c).Range\n    $cellRange.Text = $values[$c - 1]\n    $cellRange.Font.Size = 14\n    $cellRange.Font.SizeBi = 14\n}\n"}
